$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Copy the (soon to be updated) data range onto the originally-blank
# Sheet2 so the surviving sheet never picked up Sheet1's custom
# <cols> width/bestFit formatting.
[void]$ws1.Range("A1:D4").Copy()
[void]$ws2.Range("A1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Drop Sheet3 and the old Sheet1, then rename Sheet2 -> Sheet1 so the
# workbook ends up with a single sheet named "Sheet1" again.
[void]$wb.Worksheets.Item("Sheet3").Delete()
[void]$ws1.Delete()

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Sheet1"

# Update the numeric data (B:D on rows 2-4) to the new values; column A
# is unchanged.
$ws.Range("B2").Value = 0.9410862337871182
$ws.Range("C2").Value = 0.26796282082712913
$ws.Range("D2").Value = -0.20628288157466085

$ws.Range("B3").Value = 0.32722528361690628
$ws.Range("C3").Value = -0.56767435316969417
$ws.Range("D3").Value = 0.75542666256573476

$ws.Range("B4").Value = 0.085324758061261913
$ws.Range("C4").Value = -0.77842260720494083
$ws.Range("D4").Value = -0.62191472908594925

# Match the saved selection in the target file.
[void]$ws.Range("A1:D4").Select()
